# Updated run for publication
# Replace placeholder/rounded frequency values (row 2-5, columns B-X)
# with the precise recomputed frequencies from the updated run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0057929036929761
$ws.Range("C2").Value = 0.132512671976828
$ws.Range("D2").Value = 0.140477914554671
$ws.Range("E2").Value = 0.0760318609703114
$ws.Range("F2").Value = 0.000724112961622013
$ws.Range("G2").Value = 0.00217233888486604
$ws.Range("H2").Value = 0.861694424330195
$ws.Range("I2").Value = 0.0057929036929761
$ws.Range("J2").Value = 0.00506879073135409
$ws.Range("K2").Value = 0.860970311368573
$ws.Range("L2").Value = 0.0057929036929761
$ws.Range("M2").Value = 0.0246198406951484
$ws.Range("N2").Value = 0.00362056480811007
$ws.Range("O2").Value = 0.0101375814627082
$ws.Range("P2").Value = 0.827661115133961
$ws.Range("Q2").Value = 0.0101375814627082
$ws.Range("R2").Value = 0.00144822592324403
$ws.Range("S2").Value = 0.997103548153512
$ws.Range("T2").Value = 0.00506879073135409
$ws.Range("U2").Value = 0.436640115858074
$ws.Range("V2").Value = 0.0246198406951484
$ws.Range("W2").Value = 0.0963070238957277
$ws.Range("X2").Value = 0.00362056480811007

$ws.Range("B3").Value = 0.944967414916727
$ws.Range("C3").Value = 0.855901520637219
$ws.Range("D3").Value = 0.0137581462708182
$ws.Range("E3").Value = 0.876176683562636
$ws.Range("F3").Value = 0.984069514844316
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.0166545981173063
$ws.Range("I3").Value = 0.00289645184648805
$ws.Range("J3").Value = 0.990586531498914
$ws.Range("K3").Value = 0.0586531498913831
$ws.Range("L3").Value = 0.0275162925416365
$ws.Range("M3").Value = 0.00651701665459812
$ws.Range("N3").Value = 0.165097755249819
$ws.Range("O3").Value = 0.983345401882694
$ws.Range("P3").Value = 0.00651701665459812
$ws.Range("Q3").Value = 0.0369297610427227
$ws.Range("R3").Value = 0.994207096307024
$ws.Range("S3").Value = 0.00144822592324403
$ws.Range("T3").Value = 0.942070963070239
$ws.Range("U3").Value = 0.0238957277335264
$ws.Range("V3").Value = 0.0202751629254164
$ws.Range("W3").Value = 0.00941346850108617
$ws.Range("X3").Value = 0.00796524257784214

$ws.Range("B4").Value = 0.00651701665459812
$ws.Range("C4").Value = 0.00651701665459812
$ws.Range("D4").Value = 0.7697320782042
$ws.Range("E4").Value = 0.0354815351194786
$ws.Range("F4").Value = 0.00651701665459812
$ws.Range("G4").Value = 0.000724112961622013
$ws.Range("H4").Value = 0.116582186821144
$ws.Range("I4").Value = 0.990586531498914
$ws.Range("J4").Value = 0.00217233888486604
$ws.Range("K4").Value = 0.0709630702389573
$ws.Range("L4").Value = 0.00868935553946416
$ws.Range("M4").Value = 0.00941346850108617
$ws.Range("N4").Value = 0.00724112961622013
$ws.Range("O4").Value = 0.0057929036929761
$ws.Range("P4").Value = 0.162925416364953
$ws.Range("Q4").Value = 0.0224475018102824
$ws.Range("R4").Value = 0.000724112961622013
$ws.Range("S4").Value = 0.00144822592324403
$ws.Range("T4").Value = 0.00651701665459812
$ws.Range("U4").Value = 0.532947139753802
$ws.Range("V4").Value = 0.0260680666183925
$ws.Range("W4").Value = 0.8848660391021
$ws.Range("X4").Value = 0.984069514844316

$ws.Range("B5").Value = 0.0427226647356988
$ws.Range("C5").Value = 0.00506879073135409
$ws.Range("D5").Value = 0.0745836350470673
$ws.Range("E5").Value = 0.0123099203475742
$ws.Range("F5").Value = 0.00868935553946416
$ws.Range("G5").Value = 0.997103548153512
$ws.Range("H5").Value = 0.00434467776973208
$ws.Range("I5").Value = 0.000724112961622013
$ws.Range("J5").Value = 0.00144822592324403
$ws.Range("K5").Value = 0.00941346850108617
$ws.Range("L5").Value = 0.958001448225923
$ws.Range("M5").Value = 0.959449674149167
$ws.Range("N5").Value = 0.823316437364229
$ws.Range("O5").Value = 0.000724112961622013
$ws.Range("P5").Value = 0.00289645184648805
$ws.Range("Q5").Value = 0.930485155684287
$ws.Range("R5").Value = 0.00362056480811007
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0.0463432295438088
$ws.Range("U5").Value = 0.0057929036929761
$ws.Range("V5").Value = 0.929036929761043
$ws.Range("W5").Value = 0.00362056480811007
$ws.Range("X5").Value = 0
